$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "baseline" (B) and "DALI" (C) columns for rows 2-9 (lack dali and baseline)
$ws.Range("B2:C9").ClearContents()

# Update the "INR" (D) column values from 11.4 to 11.9 (rerun dali)
$ws.Range("D2:D9").Value = 11.9

# Update the selected cell in the sheet view
$ws.Range("F8").Select()
